# Ajustes a la planeación
# - Mueve la fecha de conclusión/revisión de "Vista de horarios de materias
#   (maestros)" (fila 12) de "Domingo 3 de mayo" a "Domingo 10 de mayo".
# - Actualiza la selección visible de la hoja hacia la fila de la descripción
#   del caso de uso (A20:E20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Retrasa una semana la fecha de conclusión (C12) y de revisión (D12).
$ws.Range("C12").Value = "Domingo 10 de mayo"
$ws.Range("D12").Value = "Domingo 10 de mayo"

# Actualiza la selección/posición visible de la vista de la hoja.
$ws.Range("A20:E20").Select()
